{"js": "// Replace the three-digit-number \u00f7 one-digit-number answer strings in the\n// practice table with the newly generated set of problems/answers.\n// Each entry of `replacements` is [oldText, newText] taken from the table\n// cells (e.g. \"455\u00f78=56, 7\" -> \"736\u00f77=105, 1\").\nconst replacements = [[\"455\u00f78=56, 7\", \"736\u00f77=105, 1\"], [\"610\u00f79=67, 7\", \"538\u00f79=59, 7\"], [\"324\u00f77=46, 2\", \"350\u00f72=175, 0\"], [\"997\u00f74=249, 1\", \"886\u00f75=177, 1\"], [\"421\u00f75=84, 1\", \"319\u00f77=45, 4\"], [\"981\u00f72=490, 1\", \"519\u00f74=129, 3\"], [\"213\u00f73=71, 0\", \"445\u00f75=89, 0\"], [\"512\u00f75=102, 2\", \"163\u00f75=32, 3\"], [\"549\u00f73=183, 0\", \"607\u00f75=121, 2\"], [\"325\u00f72=162, 1\", \"951\u00f73=317, 0\"], [\"441\u00f79=49, 0\", \"976\u00f79=108, 4\"], [\"175\u00f75=35, 0\", \"382\u00f76=63, 4\"], [\"823\u00f72=411, 1\", \"722\u00f77=103, 1\"], [\"611\u00f72=305, 1\", \"656\u00f78=82, 0\"], [\"755\u00f79=83, 8\", \"261\u00f74=65, 1\"], [\"321\u00f72=160, 1\", \"884\u00f73=294, 2\"], [\"645\u00f79=71, 6\", \"865\u00f74=216, 1\"], [\"504\u00f76=84, 0\", \"346\u00f77=49, 3\"], [\"220\u00f77=31, 3\", \"997\u00f72=498, 1\"], [\"894\u00f75=178, 4\", \"841\u00f75=168, 1\"], [\"387\u00f78=48, 3\", \"289\u00f77=41, 2\"], [\"311\u00f73=103, 2\", \"314\u00f77=44, 6\"], [\"471\u00f79=52, 3\", \"274\u00f73=91, 1\"], [\"941\u00f79=104, 5\", \"515\u00f74=128, 3\"], [\"637\u00f78=79, 5\", \"963\u00f72=481, 1\"]];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('455\u00f78=56, 7', '736\u00f77=105, 1'),\n    @('610\u00f79=67, 7', '538\u00f79=59, 7'),\n    @('324\u00f77=46, 2', '350\u00f72=175, 0'),\n    @('997\u00f74=249, 1', '886\u00f75=177, 1'),\n    @('421\u00f75=84, 1', '319\u00f77=45, 4'),\n    @('981\u00f72=490, 1', '519\u00f74=129, 3'),\n    @('213\u00f73=71, 0', '445\u00f75=89, 0'),\n    @('512\u00f75=102, 2', '163\u00f75=32, 3'),\n    @('549\u00f73=183, 0', '607\u00f75=121, 2'),\n    @('325\u00f72=162, 1', '951\u00f73=317, 0'),\n    @('441\u00f79=49, 0', '976\u00f79=108, 4'),\n    @('175\u00f75=35, 0', '382\u00f76=63, 4'),\n    @('823\u00f72=411, 1', '722\u00f77=103, 1'),\n    @('611\u00f72=305, 1', '656\u00f78=82, 0'),\n    @('755\u00f79=83, 8', '261\u00f74=65, 1'),\n    @('321\u00f72=160, 1', '884\u00f73=294, 2'),\n    @('645\u00f79=71, 6', '865\u00f74=216, 1'),\n    @('504\u00f76=84, 0', '346\u00f77=49, 3'),\n    @('220\u00f77=31, 3', '997\u00f72=498, 1'),\n    @('894\u00f75=178, 4', '841\u00f75=168, 1'),\n    @('387\u00f78=48, 3', '289\u00f77=41, 2'),\n    @('311\u00f73=103, 2', '314\u00f77=44, 6'),\n    @('471\u00f79=52, 3', '274\u00f73=91, 1'),\n    @('941\u00f79=104, 5', '515\u00f74=128, 3'),\n    @('637\u00f78=79, 5', '963\u00f72=481, 1'),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
